$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "HVPCB_BOM"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "HVPCB_BOM_MOUSER"

# Copy row 1 (title) and row 2 (headers) formatting+content from sheet1
$ws1.Range("A1:M2").Copy($ws2.Range("A1"))

# Move data rows 21:23 (items 19,20,21) from sheet1 to sheet2 rows 3:5, preserving format/value
$ws1.Range("A21:M23").Copy($ws2.Range("A3"))
# Re-apply the formula (Copy() only carries the cached value, not the formula)
$ws2.Range("M3:M5").Formula = "=ROUNDUP(L3*F3,2)"

# Clear the moved data out of the original sheet (keep formatting/style in place)
$ws1.Range("A21:M23").ClearContents()

# New title text + merge on the new sheet
$ws2.Range("A1:M1").Merge()
$ws2.Range("D8").Value = "a"
$ws2.Range("A1").Value = "HVPCB BOM (MOUSER)"

Write-Host "done"
